$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.392.64"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.483.48"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -1.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.479.98"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.29%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.941.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.211.22"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.12"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.505.77"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "345.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.16%  "

$ws.Range("E22").Value = "  -1.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.57%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.50"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("E26").Value = "  -3.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.58"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.91%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0863"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.58"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.40%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "435.15"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.58%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.88"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.06"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.08"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.49%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.312"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.55%  "

$ws.Range("E43").Value = "  -2.53%  "

$ws.Range("E44").Value = "  -6.40%  "

$ws.Range("E45").Value = "  -6.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "137.56"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.507"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.03%  "

$ws.Range("E49").Value = "  -1.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.98"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +22.66%  "

$ws.Range("E51").Value = "  -1.06%  "
